# Fruta / hortaliza, semanal
#
# Insert 5 new weekly price rows (Durazno: Early Majestic / Florida King,
# market date 2022-11-25 = serial 44890) at the top of the existing
# "Durazno" data block (rows 1170-1207), pushing the pre-existing rows
# down by 5 (to 1175-1212).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows, shifting rows 1170:1207 down to 1175:1212.
$ws.Rows("1170:1174").Insert()

# Shared (unchanged) identifying columns for all 5 new rows.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$prodId    = 100103
$producto  = "Frutos de hueso (carozo)"
$catId     = 100103004
$categoria = "Durazno"

$newRows = @(
    @{ Row=1170; Fecha=44890; Variedad="Early Majestic"; Calidad="Primera"; Volumen=30;  PMin=420000; PMax=430000; PProm=425000; Unidad="$/bins (420 kilos)"; Origen="Provincia de Limarí";   PKg=1012; KgUnidad=420 },
    @{ Row=1171; Fecha=44890; Variedad="Early Majestic"; Calidad="Primera"; Volumen=8;   PMin=440000; PMax=440000; PProm=440000; Unidad="$/bins (420 kilos)"; Origen="Región de O'Higgins";    PKg=1048; KgUnidad=420 },
    @{ Row=1172; Fecha=44890; Variedad="Early Majestic"; Calidad="Segunda"; Volumen=25;  PMin=400000; PMax=400000; PProm=400000; Unidad="$/bins (420 kilos)"; Origen="Provincia de Limarí";   PKg=952;  KgUnidad=420 },
    @{ Row=1173; Fecha=44890; Variedad="Early Majestic"; Calidad="Segunda"; Volumen=10;  PMin=400000; PMax=400000; PProm=400000; Unidad="$/bins (420 kilos)"; Origen="Región de O'Higgins";    PKg=952;  KgUnidad=420 },
    @{ Row=1174; Fecha=44890; Variedad="Florida King";   Calidad="Primera"; Volumen=10;  PMin=400000; PMax=400000; PProm=400000; Unidad="$/bins (420 kilos)"; Origen="Región de O'Higgins";    PKg=952;  KgUnidad=420 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
